# Resize/reposition the "Notes" textboxes (header + body) on the two
# two-column slides (slide 19 and slide 20) so they both sit at the same
# right-column position/width used elsewhere in the deck (e.g. slide 16):
#   Left  = 4648200 EMU (366.0 pt)
#   Width = 4495799 EMU (~353.99995 pt -- nudged off the exact quotient so
#           the float64->EMU round-trip lands on 4495799 instead of
#           4495798/4495800)
# Slide 20's notes-body box also gets its Height corrected from 274320 EMU
# (21.6pt, wrongly copied from the header box) to 342900 EMU (27pt), matching
# the header/body pairing used on every other "Notes" section in the deck.

$p = $ppt.ActivePresentation

$targetLeft  = 366.0
$targetWidth = 353.99995

foreach ($slideIdx in 19, 20) {
    $s = $p.Slides.Item($slideIdx)

    # "TextBox 5" = the "Notes" header label
    $notesHeader = $s.Shapes.Item(5)
    $notesHeader.Left  = $targetLeft
    $notesHeader.Width = $targetWidth

    # "TextBox 6" = the notes body paragraph
    $notesBody = $s.Shapes.Item(6)
    $notesBody.Left  = $targetLeft
    $notesBody.Width = $targetWidth
}

# Fix the notes-body height on slide 20 (was incorrectly left at the
# header's 21.6pt instead of the body's 27pt).
$slide20Body = $p.Slides.Item(20).Shapes.Item(6)
$slide20Body.Height = 27.0
